# Add Hamburger-Cord List as reference, update code and test data
$wb = $excel.ActiveWorkbook

# --- Sheet "DQ_Report": replace test data rows 2-11, drop rows 12-17 ---
$ws = $wb.Worksheets.Item("DQ_Report")

# Clear everything below the header first (old sheet used rows 2-17)
$ws.Range("A2:D17").ClearContents()

# New data rows (PatientIdentifikator, ICD_Primärkode, Orpha_Kode, dq_msg)
# Use $null for cells that must stay empty (no B or no C value)
$rows = @(
    @("P_19285751", "E84.0", 587,  "Relation  E84.0 - 587  ist im BfArM nicht vorhanden "),
    @("P_19285753", "E84.80", 587, "Relation  E84.80 - 587  ist im BfArM nicht vorhanden "),
    @("P_19285754", "E85.0", 586,  "Relation  E85.0 - 586  ist im BfArM nicht vorhanden "),
    @("P_19285755", "E75.2", 325,  "Relation  E75.2 - 325  ist im BfArM nicht vorhanden "),
    @("P_19285756", "E75.2", 320,  "Relation  E75.2 - 320  ist im BfArM nicht vorhanden "),
    @("P_19285757", $null,  586,   "Fehlendes ICD10 Code  "),
    @("P_19285758", $null,  587,   "Orpha Kodierung  587  ist im BfArM-Mapping nicht enthalten Fehlendes ICD10 Code  "),
    @("P_19285759", "E75.2", $null, "ICD10-Kodierung nicht eindeutig E75.2 Fehlendes Orpha_Kode  "),
    @("P_19285759", "E84.0", $null, "Fehlendes Orpha_Kode  "),
    @("P_19285759", "D45",   $null, "Fehlendes Orpha_Kode  ")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    if ($row[1] -ne $null) {
        $ws.Cells.Item($r, 2).Value = $row[1]
    }
    if ($row[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# --- Sheet "Statistik": update the summary row ---
$ws2 = $wb.Worksheets.Item("Statistik")
$ws2.Cells.Item(2, 2).Value = 0.28
$ws2.Cells.Item(2, 3).Value = 99.72
$ws2.Cells.Item(2, 4).Value = 90.62
$ws2.Cells.Item(2, 5).Value = 98.3
$ws2.Cells.Item(2, 6).Value = 32
